$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update VALOR MORA total (E11) and counts (C13, F13)
$ws.Range("E11").Value = 358400
$ws.Range("C13").Value = 2
$ws.Range("F13").Value = 7

# Insert 6 new rows before the existing "NOMBRE DEL REPRESENTANTE LEGAL" block (rows 21-22)
$ws.Range("A21:A26").EntireRow.Insert() | Out-Null

# New worker data rows 17-22 : CC 45466162 NELLY DEL CARMEN PEREZ MILLARES, periods 2507..2502
$periods = @("2507","2506","2505","2504","2503","2502")
for ($i = 0; $i -lt 6; $i++) {
    $r = 17 + $i
    $ws.Range("B$r").Value = "CC"
    $ws.Range("C$r").Value = "45466162"
    $ws.Range("D$r").Value = "NELLY DEL CARMEN PEREZ MILLARES"
    $ws.Range("E$r").Value = $periods[$i]
    $ws.Range("F$r").Value = 52000
    $ws.Range("G$r").Value = 1300000

    $ws.Range("B$r").Style = "Normal"
    $ws.Cells.Item($r, 2).Font.Bold = $false
}

# Copy formatting from row 16 (existing data row) to the new rows 17-22
$ws.Range("B16:J16").Copy() | Out-Null
$ws.Range("B17:J22").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Re-apply the values after paste-format (PasteSpecial formats only, values already set above) 
# Re-set values just in case paste affected them
for ($i = 0; $i -lt 6; $i++) {
    $r = 17 + $i
    $ws.Range("B$r").Value = "CC"
    $ws.Range("C$r").Value = "45466162"
    $ws.Range("D$r").Value = "NELLY DEL CARMEN PEREZ MILLARES"
    $ws.Range("E$r").Value = $periods[$i]
    $ws.Range("F$r").Value = 52000
    $ws.Range("G$r").Value = 1300000
}

$wb.Save()
